$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Steps = 0) previously held "-" / "-%" placeholder text in D3:J3.
# Replace with the actual base (steps=0) reproduction results.
$ws.Range("D3").Value = 31.0
$ws.Range("E3").Value = 0.315608888237815
$ws.Range("F3").Value = 0.670728230590925
$ws.Range("G3").Value = 0.014
$ws.Range("H3").Value = 0.899243414402008
$ws.Range("I3").Value = 0.0006828703703703704
$ws.Range("J3").Value = 9.923254

# Match the number formats already used by the other data rows for these columns.
$ws.Range("E3:H3").NumberFormat = "0.000%"
$ws.Range("I3").NumberFormat = "h:mm:ss AM/PM"
